# Apply the "changes done in existing programs" edit:
#  - Sheet1 (Billing/"CustomerName" sheet): A2 "fbgfhngfj" -> "mhgfxxfbchgj"
#  - Sheet2 (ProjectName sheet): A2 "fbgfhngfj" -> "mhgfxxfbchgj" (B2 text unchanged)
#  - Sheet3 (BillName sheet): add a yellow header-fill style to row 1, add a new
#    data row (row 2) with a new "Creditcatrd14" / "1 billing type has been
#    successfully added." / "pass" result row, add a width for the new column A
#  - Selections follow the freshly-entered cells.

$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item(1)
$ws2 = $wb.Worksheets.Item(2)
$ws3 = $wb.Worksheets.Item(3)

# --- Sheet1: CustomerName results -----------------------------------------
$ws1.Range("A2").Value = "mhgfxxfbchgj"

# --- Sheet2: ProjectName results -------------------------------------------
$ws2.Range("A2").Value = "mhgfxxfbchgj"
$ws2.Range("A2").Select()

# --- Sheet3: BillName results ------------------------------------------
# New column for the BillName header (was missing a width before).
$ws3.Columns.Item(1).ColumnWidth = 26.28515625

# Give the header row (row 1) the same yellow fill used elsewhere in the
# workbook, but keep the regular (non-bold) font.
$ws3.Range("A1:D1").Interior.Color = 65535

# New second row of data.
$ws3.Range("A2").Value = "Creditcatrd14"
$ws3.Range("B2").Value = "1 billing type has been successfully added."
$ws3.Range("C2").Value = "1 billing type has been successfully added."
$ws3.Range("D2").Value = "pass"

$ws3.Range("A2").Select()
